# Update ERD and Time-Record
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# C83: 1.33 -> 1
$ws.Range("C83").Value = 1

# B84/C84/D84: add a new "Help to edit website" time entry row
$ws.Range("B84").Value = "Help to edit website"
$ws.Range("C84").Value = 1
$ws.Range("D84").Value = "Sarvan Amel"

# Copy the number format (style) used on the row above so the new row
# matches the rest of the table.
$ws.Range("C83").Copy()
$ws.Range("C84").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("D83").Copy()
$ws.Range("D84").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# B91: extend sum formula to include C84
$ws.Range("B91").Formula = "=C75+C76+C79+C80+C81+C83+C84"

# B92 total recalculates automatically from =B89+B90+B91

# Update the view so it matches the new scroll/selection position
$excel.ActiveWindow.ScrollRow = 69
$ws.Range("B92").Select()
